# Update grades_tally and L16
# - Merge C5 ("Summarized from ") and C6 ("grades_gender.xlsx") into a single
#   cell C5 ("Summarized from grades_gender.xlsx"), and clear C6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C5").Value = "Summarized from grades_gender.xlsx"
$ws.Range("C6").ClearContents()
